$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Ají" sheet. It lands at
# row 251, pushing the previous rows 251-254 down to 252-255 (dimension
# grows from A1:R254 to A1:R255). Insert a row above the current 251 so
# the existing data shifts down intact, then fill the freshly inserted
# row with the new record's values.
$ws.Rows.Item(251).Insert()

$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 44890
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112021
$ws.Range("G251").Value = "Ají"
$ws.Range("H251").Value = "Americana (o)"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 150
$ws.Range("K251").Value = 15000
$ws.Range("L251").Value = 15000
$ws.Range("M251").Value = 15000
$ws.Range("N251").Value = "`$/caja 15 kilos"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 1000
$ws.Range("Q251").Value = 15
$ws.Range("R251").Value = "Hortaliza"
